$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price/Volume columns so numeric-looking strings
# (e.g. "29.207.93", "0.9999") are not auto-converted to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

# Coin / Link columns (plain text, unaffected by numeric coercion)
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("B10").Value = 'Solana'
$ws.Range("C10").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("B14").Value = 'Polygon'
$ws.Range("C14").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("B15").Value = 'Litecoin'
$ws.Range("C15").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("B16").Value = 'Uniswap'
$ws.Range("C16").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("B19").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C19").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("B21").Value = 'Avalanche'
$ws.Range("C21").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("B29").Value = 'PancakeSwap'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("B30").Value = 'Filecoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("B32").Value = 'Toncoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("B47").Value = 'BabyDogeCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'

# Price column (D)
$ws.Range("D2").Value = '29.207.93'
$ws.Range("D3").Value = '1.840.89'
$ws.Range("D4").Value = '0.9999'
$ws.Range("D5").Value = '242.40'
$ws.Range("D6").Value = '0.6624'
$ws.Range("D8").Value = '0.07439'
$ws.Range("D9").Value = '0.2952'
$ws.Range("D10").Value = '23.33'
$ws.Range("D11").Value = '0.07771'
$ws.Range("D12").Value = '1.880.66'
$ws.Range("D13").Value = '5.022'
$ws.Range("D14").Value = '0.6724'
$ws.Range("D15").Value = '83.48'
$ws.Range("D16").Value = '6.176'
$ws.Range("D17").Value = '0.000008562'
$ws.Range("D18").Value = '29.223.80'
$ws.Range("D19").Value = '2.097.41'
$ws.Range("D20").Value = '227.78'
$ws.Range("D21").Value = '12.55'
$ws.Range("D22").Value = '1.001'
$ws.Range("D23").Value = '7.173'
$ws.Range("D24").Value = '1.001'
$ws.Range("D25").Value = '159.19'
$ws.Range("D26").Value = '8.624'
$ws.Range("D27").Value = '0.1403'
$ws.Range("D28").Value = '18.09'
$ws.Range("D29").Value = '1.512'
$ws.Range("D30").Value = '4.136'
$ws.Range("D31").Value = '4.058'
$ws.Range("D32").Value = '1.190'
$ws.Range("D33").Value = '0.05325'
$ws.Range("D34").Value = '1.877'
$ws.Range("D35").Value = '0.7478'
$ws.Range("D36").Value = '1.154'
$ws.Range("D37").Value = '2.653'
$ws.Range("D38").Value = '1.316.55'
$ws.Range("D39").Value = '0.01800'
$ws.Range("D40").Value = '2.757'
$ws.Range("D41").Value = '6.402'
$ws.Range("D42").Value = '0.9157'
$ws.Range("D43").Value = '1.000'
$ws.Range("D44").Value = '103.05'
$ws.Range("D45").Value = '2.008.94'
$ws.Range("D46").Value = '66.01'
$ws.Range("D47").Value = '0.00000000124'
$ws.Range("D48").Value = '0.5140'
$ws.Range("D49").Value = '0.07641'
$ws.Range("D50").Value = '1.752'
$ws.Range("D51").Value = '0.05861'

# Volume(1h) column (E)
$ws.Range("E2").Value = '  +0.12%  '
$ws.Range("E3").Value = '  +0.29%  '
$ws.Range("E5").Value = '  +0.65%  '
$ws.Range("E6").Value = '  -0.70%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +0.73%  '
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("E10").Value = '  +2.13%  '
$ws.Range("E11").Value = '  +0.75%  '
$ws.Range("E12").Value = '  +4.17%  '
$ws.Range("E13").Value = '  +0.11%  '
$ws.Range("E14").Value = '  -0.49%  '
$ws.Range("E15").Value = '  -3.25%  '
$ws.Range("E16").Value = '  -0.36%  '
$ws.Range("E17").Value = '  +4.00%  '
$ws.Range("E18").Value = '  +1.12%  '
$ws.Range("E19").Value = '  +1.12%  '
$ws.Range("E20").Value = '  -0.32%  '
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("E23").Value = '  -1.27%  '
$ws.Range("E24").Value = '  +0.10%  '
$ws.Range("E25").Value = '  -0.78%  '
$ws.Range("E26").Value = '  -0.97%  '
$ws.Range("E27").Value = '  +0.18%  '
$ws.Range("E28").Value = '  +0.39%  '
$ws.Range("E29").Value = '  +0.33%  '
$ws.Range("E30").Value = '  -1.74%  '
$ws.Range("E31").Value = '  -0.53%  '
$ws.Range("E32").Value = '  -0.11%  '
$ws.Range("E33").Value = '  -0.19%  '
$ws.Range("E34").Value = '  +0.69%  '
$ws.Range("E35").Value = '  -0.40%  '
$ws.Range("E36").Value = '  +1.58%  '
$ws.Range("E37").Value = '  -0.96%  '
$ws.Range("E38").Value = '  -0.79%  '
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("E40").Value = '  +0.59%  '
$ws.Range("E41").Value = '  +7.30%  '
$ws.Range("E42").Value = '  -0.88%  '
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("E44").Value = '  -1.25%  '
$ws.Range("E45").Value = '  +2.26%  '
$ws.Range("E46").Value = '  +3.43%  '
$ws.Range("E47").Value = '  -0.44%  '
$ws.Range("E48").Value = '  -0.59%  '
$ws.Range("E49").Value = '  -6.70%  '
$ws.Range("E50").Value = '  -0.34%  '
$ws.Range("E51").Value = '  -1.40%  '

# Reset style to Normal so no explicit NumberFormat/style sticks on these cells
$ws.Range("D2:E51").Style = "Normal"
